$d = $word.ActiveDocument
$s = $d.Styles("Normal")
$s.Font.Bold = 0
